$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new testing data value into F8, referencing a new shared string
$ws.Range("F8").Value = "wrist flexion, abduction"

# Widen column F to fit the new "wrist flexion, abduction" header/data
# (mirrors the bestFit auto-sizing Excel performed in the authored workbook)
$ws.Columns.Item(6).ColumnWidth = 23.28515625

# Update the active selection to F9, matching the final workbook state
$ws.Range("F9").Select()
